$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weo_pg")

# Region column (B) changes from ITA to MEX for all data rows
$ws.Range("B2:B13").Value = "MEX"

# Update 2023 / 2030 / 2050 values (columns L, M, N) per row
# Row 2: ncap_cost, Hydropower - large-scale unit
$ws.Range("L2").Value = 2700
$ws.Range("M2").Value = 2700
$ws.Range("N2").Value = 2700

# Row 3: ncap_cost, Solar photovoltaics - Large scale unit
$ws.Range("L3").Value = 1110
$ws.Range("M3").Value = 690
$ws.Range("N3").Value = 480

# Row 4: ncap_cost, Wind offshore
$ws.Range("L4").Value = 4060
$ws.Range("M4").Value = 2760
$ws.Range("N4").Value = 1980

# Row 5: ncap_cost, Wind onshore
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1430
$ws.Range("N5").Value = 1370

# Row 6: ncap_fom, Hydropower - large-scale unit
$ws.Range("L6").Value = 70
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 65

# Row 7: ncap_fom, Solar photovoltaics - Large scale unit
$ws.Range("L7").Value = 16
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = 16

# Row 8: ncap_fom, Wind offshore
$ws.Range("L8").Value = 120
$ws.Range("M8").Value = 95
$ws.Range("N8").Value = 70

# Row 9: ncap_fom, Wind onshore
$ws.Range("L9").Value = 38
$ws.Range("M9").Value = 36
$ws.Range("N9").Value = 36

# Rows 10-13 (ncap_iled) values are unchanged.

$wb.Save()
